# Append a fresh batch of Lancers listings (fetched 2026-01-23 01:27:22) to the
# "ランサーズ" worksheet, pushing the previously-fetched rows down, and update
# the refreshed timestamp on every row plus a couple of field tweaks on the
# rows that moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2026-01-23 01:27:22"

# Insert 4 fresh rows above the old row 5 so the 4 new listings can be placed
# there while the previous row5-row8 listings shift down to row9-row12.
$ws.Rows("5:8").Insert()

# Column H ("スキル概要") grows from width 12 to 13 to fit the new tags.
$ws.Columns.Item(8).ColumnWidth = 12.17

# --- Refresh the "取得日時" timestamp for every data row (2-12) -----------
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = $timestamp
}

# --- New rows 5-8 (freshly scraped listings) -------------------------------
$ws.Cells.Item(5, 2).Value = "【フルスタックエンジニア募集】AWS構築+Pythonバックエンド開発"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5475657"
$ws.Cells.Item(5, 7).Value = 260
$ws.Cells.Item(5, 8).Value = "🔥Python ◆開発"

$ws.Cells.Item(6, 2).Value = "自動化システム"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5477084"
$ws.Cells.Item(6, 7).Value = 110
$ws.Cells.Item(6, 8).Value = "◆自動化"

$ws.Cells.Item(7, 2).Value = "Keepaの通知からAmazonで自動購入するシステムの開発依頼の仕事"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5477013"
$ws.Cells.Item(7, 7).Value = 75
$ws.Cells.Item(7, 8).Value = "◆開発"

$ws.Cells.Item(8, 2).Value = "Keepaの通知からAmazonで自動購入するシステムの開発依頼の仕事"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5477036"
$ws.Cells.Item(8, 7).Value = 75
$ws.Cells.Item(8, 8).Value = "◆開発"

# --- Row 10 (formerly row 6) gets an updated price/terms string -----------
$ws.Cells.Item(10, 4).Value = "10,000 円 ~ 20,000 円 / 募集期間 5 日、取引期間 1 日"

# --- Rebuild the hyperlinks so their refs/targets line up with the new rows
$ws.Cells.Hyperlinks.Delete()

$hyperlinkUrls = @{
    2  = "https://www.lancers.jp/work/detail/5450864"
    3  = "https://www.lancers.jp/work/detail/5460294"
    4  = "https://www.lancers.jp/work/detail/5460267"
    5  = "https://www.lancers.jp/work/detail/5475657"
    6  = "https://www.lancers.jp/work/detail/5477084"
    7  = "https://www.lancers.jp/work/detail/5477013"
    8  = "https://www.lancers.jp/work/detail/5477036"
    9  = "https://www.lancers.jp/work/detail/5476963"
    10 = "https://www.lancers.jp/work/detail/5476347"
    11 = "https://www.lancers.jp/work/detail/5476708"
    12 = "https://www.lancers.jp/work/detail/5476581"
}

for ($r = 2; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $hyperlinkUrls[$r]
    $cell.Value = $url
    $ws.Hyperlinks.Add($cell, $url)
}
